$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 14: 162. Find Peak Element ----
$ws.Range("A14").Value = "162. Find Peak Element"
$ws.Range("B14").Value = "Medium"
$ws.Range("B14").Interior.Color = $ws.Range("B4").Interior.Color
$ws.Range("C14").Value = "Binary Search"
$ws.Range("D14").Value = "Binary search but with custom check conditions. We partition based on the neighbors of the pivot value. An enum approach is the cleanest."
$ws.Hyperlinks.Add($ws.Range("E14"), "https://leetcode.com/problems/find-peak-element/solutions/1290642/intuition-behind-conditions-complete-explanation-diagram-binary-search/?envType=study-plan-v2&envId=top-interview-150 ") | Out-Null
$ws.Range("E14").Style = "Hyperlink"

# ---- Row 15: 77. Combinations ----
$ws.Range("A15").Value = "77. Combinations"
$ws.Range("B15").Value = "Medium"
$ws.Range("B15").Interior.Color = $ws.Range("B4").Interior.Color
$ws.Range("C15").Value = "Backtracking"
$ws.Hyperlinks.Add($ws.Range("E15"), "https://leetcode.com/problems/combinations/solutions/3845903/ex-amazon-explains-a-solution-with-a-video-python-javascript-java-and-c/?envType=study-plan-v2&envId=top-interview-150 ") | Out-Null
$ws.Range("E15").Style = "Hyperlink"
$ws.Range("D15").Value = "Backtracking. Use combination.size() == k for the validity condition, then perform backtracking loop."

# ---- View: move the active selection to D16 ----
$ws.Range("D16").Select()
